$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.512.12"
$ws.Range("E2").Value = "  -6.22%  "
$ws.Range("D3").Value = "3.279.33"
$ws.Range("E3").Value = "  -6.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.74%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("D9").Value = "3.270.59"
$ws.Range("E9").Value = "  -6.57%  "
$ws.Range("E10").Value = "  -11.09%  "
$ws.Range("E11").Value = "  -6.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.91%  "
$ws.Range("D14").Value = "640.12"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.34%  "
$ws.Range("D16").Value = "3.810.82"
$ws.Range("E16").Value = "  -6.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "65.456.02"
$ws.Range("E18").Value = "  -6.33%  "
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "3.287.75"
$ws.Range("E20").Value = "  -6.85%  "
$ws.Range("E21").Value = "  -8.86%  "
$ws.Range("E22").Value = "  -5.48%  "
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "107.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.31%  "
$ws.Range("E25").Value = "  -9.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.45%  "
$ws.Range("E27").Value = "  -8.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.53%  "
$ws.Range("D30").Value = "30.07"
$ws.Range("E30").Value = "  -8.38%  "
$ws.Range("E31").Value = "  -9.70%  "
$ws.Range("D33").Value = "10.98"
$ws.Range("E33").Value = "  -5.89%  "
$ws.Range("E34").Value = "  -5.66%  "
$ws.Range("D35").Value = "3.754.89"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "521.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.75%  "
$ws.Range("D39").Value = "0.0₃0725"
$ws.Range("E39").Value = "  -9.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.02%  "
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("E42").Value = "  -7.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "32.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -12.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.06%  "
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0411"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.61%  "
$ws.Range("E48").Value = "  -5.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.30%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +1.00%  "
